# Auto-generated edit script applying the Sagittarius_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 10000
$ws.Range("I76").Value = 10000
$ws.Range("K76").Value = 10000
$ws.Range("M76").Value = -9685

$ws.Range("H79").Value = 10000
$ws.Range("I79").Value = 10000
$ws.Range("K79").Value = 10000
$ws.Range("M79").Value = -8908

$ws.Range("H80").Value = 74610.60000000001
$ws.Range("J80").Value = 135594.38
$ws.Range("L80").Value = 406783.14
$ws.Range("N80").Value = -408779.14

$ws.Range("H83").Value = 74610.60000000001
$ws.Range("J83").Value = 135594.38
$ws.Range("L83").Value = 1220349.42
$ws.Range("N83").Value = -1230333.42

$ws.Range("H92").Value = 543.8889
$ws.Range("I92").Value = 496.2143
$ws.Range("K92").Value = 496.2143
$ws.Range("M92").Value = 751.7857

$ws.Range("H116").Value = 4997.3335
$ws.Range("J116").Value = 5496.25
$ws.Range("L116").Value = 5496.25
$ws.Range("N116").Value = -12380.25

$ws.Range("H132").Value = 1982.1111
$ws.Range("I132").Value = 1854.875
$ws.Range("K132").Value = 5564.625
$ws.Range("M132").Value = -3034.625

$ws.Range("H135").Value = 1282.909
$ws.Range("I135").Value = 1105.9048
$ws.Range("K135").Value = 9953.1432
$ws.Range("M135").Value = -7418.1432

$ws.Range("H137").Value = 4173.6665
$ws.Range("I137").Value = 5247.5
$ws.Range("J137").Value = 3314.6
$ws.Range("K137").Value = 15742.5
$ws.Range("L137").Value = 9943.799999999999
$ws.Range("M137").Value = -13192.5
$ws.Range("N137").Value = -15043.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8380.793
$ws.Range("I32").Value = 5597.39
$ws.Range("K32").Value = 5597.39
$ws.Range("M32").Value = -5310.39

$ws.Range("H61").Value = 7874.25
$ws.Range("I61").Value = 7500
$ws.Range("K61").Value = 7500
$ws.Range("M61").Value = -7288

$ws.Range("H102").Value = 1814.2727
$ws.Range("I102").Value = 1795.8
$ws.Range("J102").Value = 1999
$ws.Range("K102").Value = 1795.8
$ws.Range("L102").Value = 1999
$ws.Range("M102").Value = -173.8
$ws.Range("N102").Value = -5243

$ws.Range("H122").Value = 3398.0833
$ws.Range("I122").Value = 3388.7144
$ws.Range("K122").Value = 10166.1432
$ws.Range("M122").Value = -7716.143199999999

$ws.Range("H136").Value = 7874.25
$ws.Range("I136").Value = 7500
$ws.Range("K136").Value = 22500
$ws.Range("M136").Value = -19950

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()

$ws.Range("H20").Value = 2315.25
$ws.Range("I20").Value = 2504.3333
$ws.Range("K20").Value = 2504.3333
$ws.Range("M20").Value = -2257.3333

$ws.Range("H22").Value = 1583.1666
$ws.Range("I22").Value = 1499.8
$ws.Range("K22").Value = 1499.8
$ws.Range("M22").Value = -1326.8

$ws.Range("H86").Value = 4370.2856
$ws.Range("I86").Value = 2977.4
$ws.Range("J86").Value = 5144.1113
$ws.Range("K86").Value = 2977.4
$ws.Range("L86").Value = 5144.1113
$ws.Range("M86").Value = -1854.4
$ws.Range("N86").Value = -7390.1113

$ws.Range("H89").Value = 4370.2856
$ws.Range("I89").Value = 2977.4
$ws.Range("J89").Value = 5144.1113
$ws.Range("K89").Value = 14887
$ws.Range("L89").Value = 25720.5565
$ws.Range("M89").Value = -9271
$ws.Range("N89").Value = -36952.5565

$ws.Range("H105").Value = 1731.8
$ws.Range("I105").Value = 1317.125
$ws.Range("K105").Value = 1317.125
$ws.Range("M105").Value = 429.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1500
$ws.Range("I58").Value = 1500
$ws.Range("J58").Value = 1500
$ws.Range("K58").Value = 1500
$ws.Range("L58").Value = 1500
$ws.Range("M58").Value = -1297
$ws.Range("N58").Value = -1906

$ws.Range("H105").Value = 3284.6191
$ws.Range("I105").Value = 2897.7
$ws.Range("J105").Value = 3636.3635
$ws.Range("K105").Value = 2897.7
$ws.Range("L105").Value = 3636.3635
$ws.Range("M105").Value = -1150.7
$ws.Range("N105").Value = -7130.363499999999

$ws.Range("H122").Value = 2010
$ws.Range("I122").Value = 2010
$ws.Range("K122").Value = 6030
$ws.Range("M122").Value = -3580

$ws.Range("H134").Value = 2821.7727
$ws.Range("I134").Value = 2808.1428
$ws.Range("J134").Value = 3108
$ws.Range("K134").Value = 8424.428400000001
$ws.Range("L134").Value = 9324
$ws.Range("M134").Value = -5889.428400000001
$ws.Range("N134").Value = -14394

$ws.Range("H136").Value = 1500
$ws.Range("I136").Value = 1500
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 4500
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = -1950
$ws.Range("N136").Value = -9600

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3499
$ws.Range("I34").Value = 2000
$ws.Range("J34").Value = 3998.6667
$ws.Range("K34").Value = 6000
$ws.Range("L34").Value = 11996.0001
$ws.Range("M34").Value = -5916
$ws.Range("N34").Value = -12164.0001

$ws.Range("H39").Value = 6500
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 6500
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 19500
$ws.Range("N39").Value = -20088
$ws.Range("M39").ClearContents()

$ws.Range("H55").Value = 1389.9166
$ws.Range("I55").Value = 316.66666
$ws.Range("J55").Value = 1747.6666
$ws.Range("K55").Value = 949.9999799999999
$ws.Range("L55").Value = 5242.9998
$ws.Range("M55").Value = -772.9999799999999
$ws.Range("N55").Value = -5596.9998

$ws.Range("H74").Value = 15666.667
$ws.Range("J74").Value = 15666.667
$ws.Range("L74").Value = 47000.001
$ws.Range("N74").Value = -49122.001

$ws.Range("H77").Value = 15666.667
$ws.Range("J77").Value = 15666.667
$ws.Range("L77").Value = 141000.003
$ws.Range("N77").Value = -151608.003

$ws.Range("H98").Value = 2257.5386
$ws.Range("J98").Value = 2660.889
$ws.Range("L98").Value = 7982.667
$ws.Range("N98").Value = -10978.667

$ws.Range("H131").Value = 3000
$ws.Range("J131").Value = 3000
$ws.Range("L131").Value = 9000
$ws.Range("N131").Value = -19080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 89884
$ws.Range("J103").Value = 89884
$ws.Range("L103").Value = 89884
$ws.Range("N103").Value = -92228

$ws.Range("H113").Value = 1840.8572
$ws.Range("I113").Value = 1814.3334
$ws.Range("K113").Value = 1814.3334
$ws.Range("M113").Value = 355.6666

$ws.Range("H114").Value = 30000
$ws.Range("J114").Value = 30000
$ws.Range("L114").Value = 30000
$ws.Range("N114").Value = -38678

$ws.Range("H132").Value = 2723
$ws.Range("J132").Value = 1600
$ws.Range("L132").Value = 4800
$ws.Range("N132").Value = -9860

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 1000000
$ws.Range("I4").Value = 1000000
$ws.Range("K4").Value = 1000000
$ws.Range("M4").Value = -999887

$ws.Range("H16").Value = 559.6667
$ws.Range("I16").Value = 559.6667
$ws.Range("K16").Value = 559.6667
$ws.Range("M16").Value = -389.6667

$ws.Range("H28").Value = 1000000
$ws.Range("I28").Value = 1000000
$ws.Range("K28").Value = 1000000
$ws.Range("M28").Value = -999768

$ws.Range("H37").Value = 1000000
$ws.Range("I37").Value = 1000000
$ws.Range("K37").Value = 1000000
$ws.Range("M37").Value = -999893

$ws.Range("H46").Value = 3420
$ws.Range("I46").Value = 2946.5
$ws.Range("K46").Value = 2946.5
$ws.Range("M46").Value = -2758.5

$ws.Range("H122").Value = 7592.2188
$ws.Range("I122").Value = 7999.778
$ws.Range("J122").Value = 7068.2144
$ws.Range("K122").Value = 23999.334
$ws.Range("L122").Value = 21204.6432
$ws.Range("M122").Value = -21549.334
$ws.Range("N122").Value = -26104.6432

$ws.Range("H132").Value = 5115.9
$ws.Range("I132").Value = 4541.3335
$ws.Range("J132").Value = 5977.75
$ws.Range("K132").Value = 13624.0005
$ws.Range("L132").Value = 17933.25
$ws.Range("M132").Value = -11094.0005
$ws.Range("N132").Value = -22993.25

$ws.Range("H136").Value = 9296
$ws.Range("I136").Value = 10695.167
$ws.Range("K136").Value = 32085.501
$ws.Range("M136").Value = -29535.501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 97234.336
$ws.Range("J46").Value = 97234.336
$ws.Range("L46").Value = 97234.336
$ws.Range("N46").Value = -97696.336

$ws.Range("H102").Value = 80224.336
$ws.Range("J102").Value = 80224.336
$ws.Range("L102").Value = 80224.336
$ws.Range("N102").Value = -86714.336

$ws.Range("H122").Value = 3851.5
$ws.Range("I122").Value = 2451
$ws.Range("J122").Value = 4318.3335
$ws.Range("K122").Value = 7353
$ws.Range("L122").Value = 12955.0005
$ws.Range("M122").Value = -4903
$ws.Range("N122").Value = -17855.0005

$ws.Range("H132").Value = 6559.45
$ws.Range("I132").Value = 6378.3687
$ws.Range("K132").Value = 19135.1061
$ws.Range("M132").Value = -16605.1061

$ws.Range("H134").Value = 97234.336
$ws.Range("J134").Value = 97234.336
$ws.Range("L134").Value = 291703.008
$ws.Range("N134").Value = -296773.008
